$wb = $excel.ActiveWorkbook

$wsPrecond = $wb.Worksheets.Item("Precondiciones")
$wsPasos = $wb.Worksheets.Item("Pasos")

# Remove the two obsolete precondition rows (Colon address / "9" numbering cases)
$wsPrecond.Rows("3:4").Delete()

# Fix typos / wording in the "Pasos" sheet step descriptions
$wsPasos.Range("B3").Value2 = 'Ingreso "Córdoba" en el campo nombre de ciudad'
$wsPasos.Range("B5").Value2 = 'Ingreso "Colon" en el campo calle'

# Restore selections to match the saved state of each sheet
[void]$wsPrecond.Range("A3:B4").Select()
[void]$wsPasos.Select()
[void]$wsPasos.Range("B8").Select()
